$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected columns to text format so that values like
# "27.440.72", "1.004", dates, or percentages are stored as literal text
# (matching the original inlineStr cell type) rather than being
# auto-converted by Excel into numbers/dates.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.375.87'
$ws.Range('E2').Value = '  -2.81%  '
$ws.Range('D3').Value = '1.736.82'
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '323.63'
$ws.Range('E5').Value = '  -4.40%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.4227'
$ws.Range('E7').Value = '  -9.21%  '
$ws.Range('D8').Value = '0.3609'
$ws.Range('E8').Value = '  -2.97%  '
$ws.Range('D9').Value = '45.04'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.07477'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '1.123'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '21.54'
$ws.Range('E13').Value = '  -4.68%  '
$ws.Range('D14').Value = '6.076'
$ws.Range('E14').Value = '  -5.14%  '
$ws.Range('D15').Value = '7.175'
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').Value = '1.736.13'
$ws.Range('E16').Value = '  -3.45%  '
$ws.Range('D17').Value = '0.00001065'
$ws.Range('E17').Value = '  -2.87%  '
$ws.Range('D18').Value = '87.22'
$ws.Range('E18').Value = '  +5.95%  '
$ws.Range('D19').Value = '0.06045'
$ws.Range('E19').Value = '  -10.58%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '16.82'
$ws.Range('E21').Value = '  -4.13%  '
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range('D23').Value = '0.5220'
$ws.Range('E23').Value = '  -5.09%  '
$ws.Range('D24').Value = '27.427.78'
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('D25').Value = '11.33'
$ws.Range('E25').Value = '  -5.37%  '
$ws.Range('D26').Value = '2.419'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.389'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '20.17'
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('D29').Value = '149.64'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').Value = '1.935.59'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('D31').Value = '1.276'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = '126.66'
$ws.Range('E32').Value = '  -5.95%  '
$ws.Range('D33').Value = '3.732'
$ws.Range('E33').Value = '  -8.10%  '
$ws.Range('D34').Value = '5.608'
$ws.Range('E34').Value = '  -5.92%  '
$ws.Range('E35').Value = '  -6.14%  '
$ws.Range('D36').Value = '12.47'
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('D37').Value = '0.2151'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').Value = '0.06162'
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').Value = '0.02282'
$ws.Range('E39').Value = '  -4.88%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6417'
$ws.Range('E40').Value = '  -4.59%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = '5.030'
$ws.Range('E41').Value = '  -4.83%  '
$ws.Range('D42').Value = '1.188'
$ws.Range('E42').Value = '  -4.10%  '
$ws.Range('E43').Value = '  -5.00%  '
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '7.932'
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('E46').Value = '  -4.83%  '
$ws.Range('D47').Value = '3.745'
$ws.Range('D48').Value = '0.5852'
$ws.Range('E48').Value = '  -5.04%  '
$ws.Range('D49').Value = '125.46'
$ws.Range('E49').Value = '  -3.54%  '
$ws.Range('D50').Value = '1.937'
$ws.Range('E50').Value = '  -6.62%  '
$ws.Range('D51').Value = '0.06834'
$ws.Range('E51').Value = '  -4.00%  '
